$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.621.16"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "3.045.38"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'384.47"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "'102.87"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Value = "'0.544"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.590"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "'36.92"
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D13").Value = "3.518.49"
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("D14").Value = "'18.74"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").Value = "'7.76"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "3.034.75"
$ws.Range("E16").Value = "  +2.94%  "
$ws.Range("D17").Value = "'0.977"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").Value = "'10.60"
$ws.Range("E18").Value = "  -8.66%  "
$ws.Range("D19").Value = "51.668.69"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").Value = "'3.10"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").Value = "'12.44"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "'70.01"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'267.06"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'3.18"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("D26").Value = "'8.41"
$ws.Range("E26").Value = "  +6.59%  "
$ws.Range("D27").Value = "'7.39"
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("D28").Value = "'0.173"
$ws.Range("E28").Value = "  +4.27%  "
$ws.Range("D29").Value = "'26.37"
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("D32").Value = "'10.29"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'34.09"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "'2.07"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").Value = "'50.64"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("E36").Value = "  +2.29%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "'3.38"
$ws.Range("E38").Value = "  +4.04%  "
$ws.Range("E39").Value = "  +5.63%  "
$ws.Range("D40").Value = "'17.03"
$ws.Range("E40").Value = "  +2.91%  "
$ws.Range("D41").Value = "'1.87"
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("D42").Value = "'128.80"
$ws.Range("E42").Value = "  +3.29%  "
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "'2.53"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("E45").Value = "  +3.91%  "
$ws.Range("D46").Value = "'21.72"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").Value = "'2.51"
$ws.Range("E47").Value = "  +5.89%  "
$ws.Range("D48").Value = "'2.08"
$ws.Range("E48").Value = "  +3.25%  "
$ws.Range("D49").Value = "2.033.80"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").Value = "3.342.96"
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("E51").Value = "  +7.32%  "
